# Updated cryptos list on Thu Dec 21 09:22:01 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) text values for
# the crypto rows on the active worksheet.
#
# Some "Price" values are plain decimal numbers (e.g. "265.26"). Assigning
# such a string straight to Range.Value lets Excel auto-convert it into a
# numeric cell, which would lose the original text formatting. To keep the
# cell a text value (matching the source data, which stores prices as
# strings), NumberFormat is forced to "@" (Text) right before those
# assignments. Values that already contain multiple "." separators (e.g.
# "43.780.72") are never auto-parsed as numbers, so they are assigned as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.780.72"
$ws.Range("E2").Value = "  +2.12%  "

$ws.Range("D3").Value = "2.214.40"
$ws.Range("E3").Value = "  -0.12%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "265.26"
$ws.Range("E5").Value = "  +3.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.45"
$ws.Range("E6").Value = "  +12.50%  "

$ws.Range("E7").Value = "  +0.43%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("E9").Value = "  +1.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.30"
$ws.Range("E10").Value = "  +9.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0921"
$ws.Range("E11").Value = "  +1.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.62"
$ws.Range("E12").Value = "  +9.09%  "

$ws.Range("E13").Value = "  +2.10%  "

$ws.Range("D14").Value = "2.545.22"
$ws.Range("E14").Value = "  +0.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.67"
$ws.Range("E15").Value = "  +1.15%  "

$ws.Range("D16").Value = "2.201.71"
$ws.Range("E16").Value = "  -0.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.787"
$ws.Range("E17").Value = "  +0.43%  "

$ws.Range("D18").Value = "43.727.27"
$ws.Range("E18").Value = "  +2.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000103"
$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.99"
$ws.Range("E20").Value = "  +0.19%  "

$ws.Range("E22").Value = "  +7.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.16"
$ws.Range("E23").Value = "  +0.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.87"
$ws.Range("E24").Value = "  -4.74%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.56"
$ws.Range("E26").Value = "  +16.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.85"
$ws.Range("E27").Value = "  +0.74%  "

$ws.Range("E28").Value = "  +5.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.39"
$ws.Range("E29").Value = "  -8.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.25"
$ws.Range("E30").Value = "  +3.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.26"
$ws.Range("E31").Value = "  +1.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0895"
$ws.Range("E32").Value = "  +2.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.55"
$ws.Range("E33").Value = "  +0.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.41"
$ws.Range("E34").Value = "  +3.52%  "

$ws.Range("E35").Value = "  +1.49%  "

$ws.Range("E36").Value = "  +1.95%  "

$ws.Range("E37").Value = "  -0.94%  "

$ws.Range("E38").Value = "  +1.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.29"
$ws.Range("E39").Value = "  +17.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.37"
$ws.Range("E40").Value = "  -3.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "65.02"
$ws.Range("E41").Value = "  +8.29%  "

$ws.Range("E42").Value = "  -1.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.50"
$ws.Range("E43").Value = "  +3.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.203"
$ws.Range("E44").Value = "  +1.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.63"
$ws.Range("E45").Value = "  -2.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.35"
$ws.Range("E47").Value = "  -0.11%  "

$ws.Range("E48").Value = "  +4.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.13"
$ws.Range("E49").Value = "  +0.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.444"
$ws.Range("E50").Value = "  -4.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.50"
$ws.Range("E51").Value = "  +5.30%  "
